$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap match-detail columns (F:V) between row 91 and row 92 ---
# (Index, pais, torneio, temporada and data_partida in A:E stay the same;
# only the match info from column F onward was reordered.)
$range91 = $ws.Range("F91:V91")
$range92 = $ws.Range("F92:V92")

$v91 = $range91.Value()
$v92 = $range92.Value()

$range91.Value = $v92
$range92.Value = $v91

# --- 2) Append two new rows (132 and 133) with the same look & feel as row 131 ---
$ws.Range("A131:V131").Copy($ws.Range("A132:V132"))
$ws.Range("A131:V131").Copy($ws.Range("A133:V133"))

# Row 132
$ws.Cells.Item(132, 1).Value = 131
$ws.Cells.Item(132, 2).Value = "costa-rica"
$ws.Cells.Item(132, 3).Value = "primera-division"
$ws.Cells.Item(132, 4).Value = "2023-2024"
$ws.Cells.Item(132, 5).Value = 45256.91666666666
$ws.Cells.Item(132, 6).Value = "AD Santos"
$ws.Cells.Item(132, 7).Value = 1
$ws.Cells.Item(132, 8).Value = "Saprissa"
$ws.Cells.Item(132, 9).Value = 2
$ws.Cells.Item(132, 10).Value = 4.2
$ws.Cells.Item(132, 11).Value = "23/11/2023 22:12"
$ws.Cells.Item(132, 12).Value = 3.51
$ws.Cells.Item(132, 13).Value = "26/11/2023 21:59"
$ws.Cells.Item(132, 14).Value = 3.68
$ws.Cells.Item(132, 15).Value = "23/11/2023 22:12"
$ws.Cells.Item(132, 16).Value = 3.33
$ws.Cells.Item(132, 17).Value = "26/11/2023 21:59"
$ws.Cells.Item(132, 18).Value = 1.83
$ws.Cells.Item(132, 19).Value = "23/11/2023 22:12"
$ws.Cells.Item(132, 20).Value = 2.18
$ws.Cells.Item(132, 21).Value = "26/11/2023 21:59"
$ws.Cells.Item(132, 22).Value = "https://www.betexplorer.com/football/costa-rica/primera-division/santos-de-guapiles-saprissa/UPms3Ig2/"

# Row 133
$ws.Cells.Item(133, 1).Value = 132
$ws.Cells.Item(133, 2).Value = "costa-rica"
$ws.Cells.Item(133, 3).Value = "primera-division"
$ws.Cells.Item(133, 4).Value = "2023-2024"
$ws.Cells.Item(133, 5).Value = 45257.04166666666
$ws.Cells.Item(133, 6).Value = "Liberia"
$ws.Cells.Item(133, 7).Value = 0
$ws.Cells.Item(133, 8).Value = "San Carlos"
$ws.Cells.Item(133, 9).Value = 1
$ws.Cells.Item(133, 10).Value = 2.97
$ws.Cells.Item(133, 11).Value = "24/11/2023 02:12"
$ws.Cells.Item(133, 12).Value = 2.55
$ws.Cells.Item(133, 13).Value = "27/11/2023 00:50"
$ws.Cells.Item(133, 14).Value = 3.58
$ws.Cells.Item(133, 15).Value = "24/11/2023 02:12"
$ws.Cells.Item(133, 16).Value = 3.77
$ws.Cells.Item(133, 17).Value = "27/11/2023 00:50"
$ws.Cells.Item(133, 18).Value = 2.19
$ws.Cells.Item(133, 19).Value = "24/11/2023 02:12"
$ws.Cells.Item(133, 20).Value = 2.58
$ws.Cells.Item(133, 21).Value = "27/11/2023 00:50"
$ws.Cells.Item(133, 22).Value = "https://www.betexplorer.com/football/costa-rica/primera-division/liberia-san-carlos/xSUXPu1F/"

$wb.Save()
